# Add an "asset_type" column (F) to the portfolio data sheet, classifying
# each holding as either "Stock" or "ETF" based on its existing sector.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting from E1 (bold, centered, bordered) onto F1,
# then set the header text.
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "asset_type"

$ws.Range("F2").Value = "Stock"
$ws.Range("F3").Value = "Stock"
$ws.Range("F4").Value = "Stock"
$ws.Range("F5").Value = "ETF"
$ws.Range("F6").Value = "ETF"
$ws.Range("F7").Value = "ETF"
$ws.Range("F8").Value = "Stock"
$ws.Range("F9").Value = "Stock"
$ws.Range("F10").Value = "Stock"
$ws.Range("F11").Value = "Stock"
$ws.Range("F12").Value = "ETF"
$ws.Range("F13").Value = "ETF"
$ws.Range("F14").Value = "Stock"
$ws.Range("F15").Value = "Stock"
$ws.Range("F16").Value = "Stock"
$ws.Range("F17").Value = "Stock"
